$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new rows before row 8, pushing the existing data (rows 8-30) down to rows 11-33
$ws.Rows("8:10").Insert()

# Fill the 3 new rows (8, 9, 10) with the new weekly data
$ws.Range("A8").Value = 12
$ws.Range("B8").Value = "Mapocho Venta Directa de Santiago"
$ws.Range("C8").Value = "Metropolitana"
$ws.Range("D8").Value = 44467
$ws.Range("E8").Value = 13
$ws.Range("F8").Value = 300000000
$ws.Range("G8").Value = "Espárragos"
$ws.Range("H8").Value = "Sin especificar"
$ws.Range("I8").Value = "Banquete"
$ws.Range("J8").Value = 450
$ws.Range("K8").Value = 1600
$ws.Range("L8").Value = 1600
$ws.Range("M8").Value = 1600
$ws.Range("N8").Value = "$/kilo"
$ws.Range("O8").Value = "Provincia de Linares"
$ws.Range("P8").Value = 1600
$ws.Range("Q8").Value = 1
$ws.Range("R8").Value = "Hortaliza"

$ws.Range("A9").Value = 12
$ws.Range("B9").Value = "Mapocho Venta Directa de Santiago"
$ws.Range("C9").Value = "Metropolitana"
$ws.Range("D9").Value = 44467
$ws.Range("E9").Value = 13
$ws.Range("F9").Value = 300000000
$ws.Range("G9").Value = "Espárragos"
$ws.Range("H9").Value = "Sin especificar"
$ws.Range("I9").Value = "Primera"
$ws.Range("J9").Value = 440
$ws.Range("K9").Value = 1500
$ws.Range("L9").Value = 1500
$ws.Range("M9").Value = 1500
$ws.Range("N9").Value = "$/kilo"
$ws.Range("O9").Value = "Provincia de Linares"
$ws.Range("P9").Value = 1500
$ws.Range("Q9").Value = 1
$ws.Range("R9").Value = "Hortaliza"

$ws.Range("A10").Value = 12
$ws.Range("B10").Value = "Mapocho Venta Directa de Santiago"
$ws.Range("C10").Value = "Metropolitana"
$ws.Range("D10").Value = 44467
$ws.Range("E10").Value = 13
$ws.Range("F10").Value = 300000000
$ws.Range("G10").Value = "Espárragos"
$ws.Range("H10").Value = "Sin especificar"
$ws.Range("I10").Value = "Segunda"
$ws.Range("J10").Value = 380
$ws.Range("K10").Value = 1300
$ws.Range("L10").Value = 1300
$ws.Range("M10").Value = 1300
$ws.Range("N10").Value = "$/kilo"
$ws.Range("O10").Value = "Provincia de Linares"
$ws.Range("P10").Value = 1300
$ws.Range("Q10").Value = 1
$ws.Range("R10").Value = "Hortaliza"
